# Scheduled runner update: refresh market-price-derived profit columns (H:N)
# across the Leve profit tables. Values sourced from the latest price snapshot.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 6 (Leve Item ID 4564)
$ws.Range("H6").Value = 1499.75
$ws.Range("I6").Value = 1866.3334
$ws.Range("J6").Value = 400
$ws.Range("K6").Value = 5599.0002
$ws.Range("L6").Value = 1200
$ws.Range("M6").Value = -5487.0002
$ws.Range("N6").Value = -1424
# Row 39 (Leve Item ID 4603)
$ws.Range("H39").Value = 787.5
$ws.Range("I39").Value = 575
$ws.Range("J39").Value = 1000
$ws.Range("K39").Value = 1725
$ws.Range("L39").Value = 3000
$ws.Range("M39").Value = -1429
$ws.Range("N39").Value = -3592
# Row 86 (Leve Item ID 12603)
$ws.Range("H86").Value = 4053.182
$ws.Range("I86").Value = 3975
$ws.Range("J86").Value = 4097.857
$ws.Range("K86").Value = 3975
$ws.Range("L86").Value = 4097.857
$ws.Range("M86").Value = -2852
$ws.Range("N86").Value = -6343.857
# Row 89 (Leve Item ID 12603)
$ws.Range("H89").Value = 4053.182
$ws.Range("I89").Value = 3975
$ws.Range("J89").Value = 4097.857
$ws.Range("K89").Value = 19875
$ws.Range("L89").Value = 20489.285
$ws.Range("M89").Value = -14259
$ws.Range("N89").Value = -31721.285
# Row 93 (Leve Item ID 18043)
$ws.Range("H93").Value = 26601
$ws.Range("J93").Value = 26601
$ws.Range("L93").Value = 26601
$ws.Range("N93").Value = -31593
# Row 101 (Leve Item ID 19884)
$ws.Range("H101").Value = 584
$ws.Range("I101").Value = 482.5
$ws.Range("K101").Value = 1447.5
$ws.Range("M101").Value = 174.5
# Row 103 (Leve Item ID 19909)
$ws.Range("H103").Value = 2012.8334
$ws.Range("I103").Value = 1539
$ws.Range("J103").Value = 2249.75
$ws.Range("K103").Value = 4617
$ws.Range("L103").Value = 6749.25
$ws.Range("M103").Value = -4031
$ws.Range("N103").Value = -7921.25
# Row 107 (Leve Item ID 27766)
$ws.Range("H107").Value = 397.16666
$ws.Range("I107").Value = 309.4375
$ws.Range("K107").Value = 309.4375
$ws.Range("M107").Value = 1610.5625
# Row 113 (Leve Item ID 27775)
$ws.Range("H113").Value = 6135
$ws.Range("I113").Value = 5905
$ws.Range("J113").Value = 6250
$ws.Range("K113").Value = 5905
$ws.Range("L113").Value = 6250
$ws.Range("M113").Value = -2651
$ws.Range("N113").Value = -12758
# Row 137 (Leve Item ID 44013)
$ws.Range("H137").Value = 2181.875
$ws.Range("I137").Value = 2109.8333
$ws.Range("K137").Value = 6329.499899999999
$ws.Range("M137").Value = -3779.499899999999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 37 (Leve Item ID 3096)
$ws.Range("H37").Value = 24245
$ws.Range("I37").Value = 3500
$ws.Range("K37").Value = 3500
$ws.Range("M37").Value = -3227
# Row 61 (Leve Item ID 43999)
$ws.Range("H61").Value = 3428
$ws.Range("I61").Value = 3237.3333
$ws.Range("J61").Value = 4000
$ws.Range("K61").Value = 3237.3333
$ws.Range("L61").Value = 4000
$ws.Range("M61").Value = -3025.3333
$ws.Range("N61").Value = -4424
# Row 74 (Leve Item ID 44000)
$ws.Range("H74").Value = 12999.6
$ws.Range("I74").Value = 15999.75
$ws.Range("J74").Value = 999
$ws.Range("K74").Value = 15999.75
$ws.Range("L74").Value = 999
$ws.Range("M74").Value = -15125.75
$ws.Range("N74").Value = -2747
# Row 77 (Leve Item ID 44000)
$ws.Range("H77").Value = 12999.6
$ws.Range("I77").Value = 15999.75
$ws.Range("J77").Value = 999
$ws.Range("K77").Value = 79998.75
$ws.Range("L77").Value = 4995
$ws.Range("M77").Value = -75630.75
$ws.Range("N77").Value = -13731
# Row 80 (Leve Item ID 10667)
$ws.Range("H80").Value = 54997.5
$ws.Range("J80").Value = 54997.5
$ws.Range("L80").Value = 54997.5
$ws.Range("N80").Value = -56993.5
# Row 83 (Leve Item ID 10667)
$ws.Range("H83").Value = 54997.5
$ws.Range("J83").Value = 54997.5
$ws.Range("L83").Value = 164992.5
$ws.Range("N83").Value = -174976.5
# Row 130 (Leve Item ID 34732)
$ws.Range("H130").Value = 96500
$ws.Range("J130").Value = 96500
$ws.Range("L130").Value = 96500
$ws.Range("N130").Value = -106540
# Row 132 (Leve Item ID 43997)
$ws.Range("H132").Value = 4896.8
$ws.Range("I132").Value = 4743.5
$ws.Range("K132").Value = 14230.5
$ws.Range("M132").Value = -11700.5
# Row 136 (Leve Item ID 43999)
$ws.Range("H136").Value = 3428
$ws.Range("I136").Value = 3237.3333
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 9711.999899999999
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = -7161.999899999999
$ws.Range("N136").Value = -17100

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 58 (Leve Item ID 44021)
$ws.Range("H58").Value = 4264.4443
$ws.Range("I58").Value = 2845
$ws.Range("J58").Value = 5400
$ws.Range("K58").Value = 2845
$ws.Range("L58").Value = 5400
$ws.Range("M58").Value = -2642
$ws.Range("N58").Value = -5806
# Row 68 (Leve Item ID 10611)
$ws.Range("H68").Value = 69765
$ws.Range("J68").Value = 69765
$ws.Range("L68").Value = 69765
$ws.Range("N68").Value = -71263
# Row 71 (Leve Item ID 10611)
$ws.Range("H71").Value = 69765
$ws.Range("J71").Value = 69765
$ws.Range("L71").Value = 209295
$ws.Range("N71").Value = -216783
# Row 103 (Leve Item ID 19558)
$ws.Range("H103").Value = 50262
$ws.Range("I103").Value = 50262
$ws.Range("K103").Value = 50262
$ws.Range("M103").Value = -49090
# Row 136 (Leve Item ID 44021)
$ws.Range("H136").Value = 4264.4443
$ws.Range("I136").Value = 2845
$ws.Range("J136").Value = 5400
$ws.Range("K136").Value = 8535
$ws.Range("L136").Value = 16200
$ws.Range("M136").Value = -5985
$ws.Range("N136").Value = -21300
# Row 141 (Leve Item ID 43345)
$ws.Range("H141").Value = 424491.25
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 424491.25
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 424491.25
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -434851.25

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 34 (Leve Item ID 4749)
$ws.Range("H34").Value = 7414.1665
$ws.Range("J34").Value = 7414.1665
$ws.Range("L34").Value = 22242.4995
$ws.Range("N34").Value = -22410.4995
# Row 39 (Leve Item ID 4712)
$ws.Range("H39").Value = 16498.084
$ws.Range("J39").Value = 16498.084
$ws.Range("L39").Value = 49494.25199999999
$ws.Range("N39").Value = -50082.25199999999
# Row 55 (Leve Item ID 4733)
$ws.Range("H55").Value = 12767.154
$ws.Range("J55").Value = 13414.417
$ws.Range("L55").Value = 40243.251
$ws.Range("N55").Value = -40597.251
# Row 118 (Leve Item ID 27872)
$ws.Range("H118").Value = 4904
$ws.Range("I118").Value = 3944
$ws.Range("K118").Value = 11832
$ws.Range("M118").Value = -10589

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 21 (Leve Item ID 4430)
$ws.Range("H21").Value = 4999.5
$ws.Range("J21").Value = 4999.5
$ws.Range("L21").Value = 4999.5
$ws.Range("N21").Value = -5345.5
# Row 30 (Leve Item ID 4430)
$ws.Range("H30").Value = 4999.5
$ws.Range("J30").Value = 4999.5
$ws.Range("L30").Value = 4999.5
$ws.Range("N30").Value = -5209.5
# Row 80 (Leve Item ID 12521)
$ws.Range("H80").Value = 4000
$ws.Range("I80").Value = 4000
$ws.Range("K80").Value = 4000
$ws.Range("M80").Value = -3002
# Row 83 (Leve Item ID 12521)
$ws.Range("H83").Value = 4000
$ws.Range("I83").Value = 4000
$ws.Range("K83").Value = 20000
$ws.Range("M83").Value = -15008
# Row 132 (Leve Item ID 44008)
$ws.Range("H132").Value = 5474.625
$ws.Range("I132").Value = 3600
$ws.Range("J132").Value = 6599.4
$ws.Range("K132").Value = 10800
$ws.Range("L132").Value = 19798.2
$ws.Range("M132").Value = -8270
$ws.Range("N132").Value = -24858.2

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 5 (Leve Item ID 3790)
$ws.Range("H5").Value = 50000
$ws.Range("J5").Value = 50000
$ws.Range("L5").Value = 50000
$ws.Range("N5").Value = -50226
# Row 24 (Leve Item ID 3774)
$ws.Range("H24").Value = 20288.375
$ws.Range("I24").Value = 15383.333
$ws.Range("J24").Value = 35003.5
$ws.Range("K24").Value = 15383.333
$ws.Range("L24").Value = 35003.5
$ws.Range("M24").Value = -15040.333
$ws.Range("N24").Value = -35689.5
# Row 132 (Leve Item ID 44058)
$ws.Range("H132").Value = 3250.3333
$ws.Range("I132").Value = 3014.2856
$ws.Range("J132").Value = 3580.8
$ws.Range("K132").Value = 9042.856800000001
$ws.Range("L132").Value = 10742.4
$ws.Range("M132").Value = -6512.856800000001
$ws.Range("N132").Value = -15802.4

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 26 (Leve Item ID 3800)
$ws.Range("H26").Value = 60000
$ws.Range("J26").Value = 60000
$ws.Range("L26").Value = 60000
$ws.Range("N26").Value = -60586
# Row 28 (Leve Item ID 3053)
$ws.Range("H28").Value = 10000
$ws.Range("I28").Value = 10000
$ws.Range("K28").Value = 10000
$ws.Range("M28").Value = -9652
# Row 75 (Leve Item ID 11957)
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("M75").ClearContents()
$ws.Range("N75").ClearContents()
# Row 78 (Leve Item ID 11957)
$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("M78").ClearContents()
$ws.Range("N78").ClearContents()
# Row 107 (Leve Item ID 27746)
$ws.Range("H107").Value = 732.5
$ws.Range("I107").Value = 615.3333
$ws.Range("J107").Value = 849.6667
$ws.Range("K107").Value = 1845.9999
$ws.Range("L107").Value = 2549.0001
$ws.Range("M107").Value = 74.00009999999997
$ws.Range("N107").Value = -6389.0001

